$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 7 entirely (it only contained "Descriptions (optional)" in A7),
# which shifts rows 8-10 up to become rows 7-9.
$ws.Rows.Item(7).Delete()

# Update the selection to match the post-edit state (merged cell B5:F5 selected).
$ws.Range("B5:F5").Select()
